$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Förändrad" (column C) for rows 2-9: 46062 -> 46063 ---
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value2 = 46063
}

# --- Reorder the data rows 4-9 (columns A, B, F, G) ---
# New row 4: previously row 8 (A 35734-2023)
$ws.Range("A4").Value2 = "A 35734-2023"
$ws.Range("B4").Value2 = 45147.89258101852
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value2 = 5.9

# New row 5: previously row 7 (A 25251-2025)
$ws.Range("A5").Value2 = "A 25251-2025"
$ws.Range("B5").Value2 = 45800.50082175926
$ws.Range("F5").Value2 = "Kommuner"
$ws.Range("G5").Value2 = 0.7

# New row 6: previously row 9 (A 5402-2026)
$ws.Range("A6").Value2 = "A 5402-2026"
$ws.Range("B6").Value2 = 46050.49721064815
$ws.Range("F6").Value2 = "Kommuner"
$ws.Range("G6").Value2 = 0.7

# New row 7: previously row 4 (A 6983-2023)
$ws.Range("A7").Value2 = "A 6983-2023"
$ws.Range("B7").Value2 = 44967.68585648148
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value2 = 5.4

# New row 8: previously row 5 (A 25254-2025)
$ws.Range("A8").Value2 = "A 25254-2025"
$ws.Range("B8").Value2 = 45800.50479166667
$ws.Range("F8").Value2 = "Kommuner"
$ws.Range("G8").Value2 = 0.2

# New row 9: previously row 6 (A 26074-2025)
$ws.Range("A9").Value2 = "A 26074-2025"
$ws.Range("B9").Value2 = 45805.32366898148
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value2 = 1.3
